$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -9
$ws.Range("F5").Value = -7
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -7
$ws.Range("F9").Value = 4
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = 2
